$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting "Localización" (D) and "Tipo" (E)
# one column to the right (Localización -> E, Tipo -> F).
$ws.Range("D1").EntireColumn.Insert()

# New column D becomes "Latitud"
$ws.Range("D1").Value = "Latitud"
$ws.Range("D2").Value = 15.56
$ws.Range("D3").Value = 15.56

# Old "Localización" column (now shifted to E) becomes "Longitud" with numeric data
$ws.Range("E1").Value = "Longitud"
$ws.Range("E2").Value = 25.26
$ws.Range("E3").Value = 25.26

# Update the used range / selection to match the new layout
$ws.Range("F2").Select()
